# Regenerate the localization-status report:
#  - two files that were previously "Ready for handoff" (017e9cf1..., 954cdb8d...)
#    have moved on to "In Translation", with a refreshed handoff timestamp
#  - two brand-new files (5e6bfec9..., a409dc41...) have just been queued and
#    show up as freshly "Ready for handoff"
#
# NOTE: in this COM-interop runtime, Range.Hyperlinks.Delete() clears every
# hyperlink on the *whole sheet* (not just the range it was called on), and
# Hyperlink.Delete() on an individual hyperlink object is a no-op. So each
# sheet's hyperlinks are cleared once up front and then ALL of them (existing
# + new) are re-added from scratch in top-to-bottom order.

$wb = $excel.ActiveWorkbook

$mdBase = "https://github.com/OpenLocalizationTest/oltest/blob/a84c75f7e233bb13e1d62060cc3a303710133f1a/e2e/"
$cfgUrl = "https://github.com/OpenLocalizationTest/oltest/blob/a84c75f7e233bb13e1d62060cc3a303710133f1a/.localization-config"
$zhBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d080868d220c0f24f6ea357ab01df54b772eb329/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/"
$deBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b63289999a5a85791d62304dbef0ecb23f9059ec/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/"

$file1 = "017e9cf1-ce6c-442c-ab1a-675b7eca1942"
$file2 = "954cdb8d-4b1c-47b1-8fd0-3a61e8cd4c98"
$file3 = "5e6bfec9-9c10-4ef5-988b-8ebc7fca0e45"
$file4 = "a409dc41-8675-45ae-bb8a-3121e5df4d6a"

$hash1 = "4a9554b03f30e35926f025ff3993439a07f0612d"
$hash2 = "9f38f52f9d5162ffde05abdbf4d42a542a1eb29b"
$hash3 = "0524561e47ee83286f38ba82c340e07d7e4a8834"
$hash4 = "297c87bf6ed6d1a1d6feac22653dff4ead02c465"

$file1Md = "$file1.md"
$file2Md = "$file2.md"
$file3Md = "$file3.md"
$file4Md = "$file4.md"

$newHandoffZh = "2016-01-13 11:23:18"
$newHandoffDe = "2016-01-13 11:23:56"

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$ovw = $wb.Worksheets.Item("Overview")

$ovw.Range("B2").Value = "In Translation"
$ovw.Range("C2").Value = "In Translation"
$ovw.Range("B3").Value = "In Translation"
$ovw.Range("C3").Value = "In Translation"

# row 4 used to hold ".localization-config" - it now holds the first newly
# queued file, and ".localization-config" shifts down to row 6
$ovw.Range("A4").Value = $file3Md
$ovw.Range("B4").Value = "Ready for handoff"
$ovw.Range("C4").Value = "Ready for handoff"

$ovw.Range("A5").Value = $file4Md
$ovw.Range("B5").Value = "Ready for handoff"
$ovw.Range("C5").Value = "Ready for handoff"

$ovw.Range("A6").Value = ".localization-config"
$ovw.Range("B6").Value = "Not to be localized"
$ovw.Range("C6").Value = "Not to be localized"

$ovw.Range("A1").Hyperlinks.Delete()
$ovw.Hyperlinks.Add($ovw.Range("A2"), "$mdBase$file1Md", "", "", $file1Md)
$ovw.Hyperlinks.Add($ovw.Range("A3"), "$mdBase$file2Md", "", "", $file2Md)
$ovw.Hyperlinks.Add($ovw.Range("A4"), "$mdBase$file3Md", "", "", $file3Md)
$ovw.Hyperlinks.Add($ovw.Range("A5"), "$mdBase$file4Md", "", "", $file4Md)
$ovw.Hyperlinks.Add($ovw.Range("A6"), $cfgUrl, "", "", ".localization-config")
$ovw.Range("A2:A6").Style = "HyperLink"

# ---------------------------------------------------------------------------
# Helper applied to each language detail sheet (zh-cn / de-de)
# ---------------------------------------------------------------------------
function Update-DetailSheet {
    param($sheet, $handoffDate, $xlfBase, $xlf1, $xlf2, $xlf3, $xlf4)

    # Row 2 / Row 3: existing files move to "In Translation" with the
    # refreshed handoff datetime
    $sheet.Range("B2").Value = "In Translation"
    $sheet.Range("D2").Value = $handoffDate
    $sheet.Range("B3").Value = "In Translation"
    $sheet.Range("D3").Value = $handoffDate

    # Row 4: new file (was the lone ".localization-config" row before)
    $sheet.Range("A4").Value = $file3Md
    $sheet.Range("B4").Value = "Ready for handoff"
    $sheet.Range("C4").Value = $xlf3
    $sheet.Range("D4").Value = $handoffDate
    $sheet.Range("G4").Value = "0001-01-01 00:00:00"
    $sheet.Range("H4").Value = "Include"

    # Row 5: second new file
    $sheet.Range("A5").Value = $file4Md
    $sheet.Range("B5").Value = "Ready for handoff"
    $sheet.Range("C5").Value = $xlf4
    $sheet.Range("D5").Value = $handoffDate
    $sheet.Range("G5").Value = "0001-01-01 00:00:00"
    $sheet.Range("H5").Value = "Include"

    # Row 6: the ".localization-config" row, pushed down from row 4
    $sheet.Range("A6").Value = ".localization-config"
    $sheet.Range("B6").Value = "Not to be localized"
    $sheet.Range("D6").Value = "0001-01-01 00:00:00"
    $sheet.Range("G6").Value = "0001-01-01 00:00:00"
    $sheet.Range("H6").Value = "Ignored"

    # Rebuild every hyperlink on the sheet from scratch (see note up top)
    $sheet.Range("A1").Hyperlinks.Delete()
    $sheet.Hyperlinks.Add($sheet.Range("A2"), "$mdBase$file1Md", "", "", $file1Md)
    $sheet.Hyperlinks.Add($sheet.Range("C2"), "$xlfBase$xlf1", "", "", $xlf1)
    $sheet.Hyperlinks.Add($sheet.Range("A3"), "$mdBase$file2Md", "", "", $file2Md)
    $sheet.Hyperlinks.Add($sheet.Range("C3"), "$xlfBase$xlf2", "", "", $xlf2)
    $sheet.Hyperlinks.Add($sheet.Range("A4"), "$mdBase$file3Md", "", "", $file3Md)
    $sheet.Hyperlinks.Add($sheet.Range("C4"), "$xlfBase$xlf3", "", "", $xlf3)
    $sheet.Hyperlinks.Add($sheet.Range("A5"), "$mdBase$file4Md", "", "", $file4Md)
    $sheet.Hyperlinks.Add($sheet.Range("C5"), "$xlfBase$xlf4", "", "", $xlf4)
    $sheet.Hyperlinks.Add($sheet.Range("A6"), $cfgUrl, "", "", ".localization-config")
    $sheet.Range("A2").Style = "HyperLink"
    $sheet.Range("C2").Style = "HyperLink"
    $sheet.Range("A3").Style = "HyperLink"
    $sheet.Range("C3").Style = "HyperLink"
    $sheet.Range("A4").Style = "HyperLink"
    $sheet.Range("C4").Style = "HyperLink"
    $sheet.Range("A5").Style = "HyperLink"
    $sheet.Range("C5").Style = "HyperLink"
    $sheet.Range("A6").Style = "HyperLink"
}

$zh = $wb.Worksheets.Item("zh-cn")
$xlf1Zh = "$file1.$hash1.zh-cn.xlf"
$xlf2Zh = "$file2.$hash2.zh-cn.xlf"
$xlf3Zh = "$file3.$hash3.zh-cn.xlf"
$xlf4Zh = "$file4.$hash4.zh-cn.xlf"
Update-DetailSheet $zh $newHandoffZh $zhBase $xlf1Zh $xlf2Zh $xlf3Zh $xlf4Zh

$de = $wb.Worksheets.Item("de-de")
$xlf1De = "$file1.$hash1.de-de.xlf"
$xlf2De = "$file2.$hash2.de-de.xlf"
$xlf3De = "$file3.$hash3.de-de.xlf"
$xlf4De = "$file4.$hash4.de-de.xlf"
Update-DetailSheet $de $newHandoffDe $deBase $xlf1De $xlf2De $xlf3De $xlf4De
